# Insert a new data row before the current row 230.
# This shifts the existing rows 230-279 down to 231-280 (preserving their
# content/formatting), and row 230 receives new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole row at 230; Excel shifts rows 230..279 down to 231..280.
$ws.Rows("230:230").Insert()

# Populate the newly inserted (blank) row 230 with the new record.
$ws.Cells.Item(230, 1).Value  = 10
$ws.Cells.Item(230, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(230, 3).Value  = "La Araucanía"
$ws.Cells.Item(230, 4).Value  = 44798
$ws.Cells.Item(230, 5).Value  = 9
$ws.Cells.Item(230, 6).Value  = 100112039
$ws.Cells.Item(230, 7).Value  = "Ciboulette"
$ws.Cells.Item(230, 8).Value  = "Sin especificar"
$ws.Cells.Item(230, 9).Value  = "Primera"
$ws.Cells.Item(230, 10).Value = 50
$ws.Cells.Item(230, 11).Value = 7000
$ws.Cells.Item(230, 12).Value = 7000
$ws.Cells.Item(230, 13).Value = 7000
$ws.Cells.Item(230, 14).Value = "$/docena de atados"
$ws.Cells.Item(230, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(230, 16).Value = 2333
$ws.Cells.Item(230, 17).Value = 3
$ws.Cells.Item(230, 18).Value = "Hortaliza"

# Apply the same date number format used by the other rows in column D.
$ws.Cells.Item(230, 4).NumberFormat = $ws.Cells.Item(231, 4).NumberFormat
